$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EPV")

# Rows 24-32 had column C accidentally duplicating the tag/label text that
# belongs in column B (e.g. "Premises and Equipment", "Current Year Revenue",
# etc.). Retrieve/clear those stray tags from column C while leaving the
# cell formatting (style) and column B labels untouched.
$ws.Range("C24:C32").ClearContents()
